# Crypto price/volume refresh - Sat Jun  3 03:30:37 UTC 2023 (GitHub Actions)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.171.26'
$ws.Range('E2').Value = '  +0.91%  '
$ws.Range('D3').Value = '1.901.70'
$ws.Range('E3').Value = '  +1.40%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').Value = '''306.81'
$ws.Range('E5').Value = '  +0.14%  '
$ws.Range('E6').Value = '  +0.18%  '
$ws.Range('D7').Value = '''0.5232'
$ws.Range('E7').Value = '  +1.39%  '
$ws.Range('E8').Value = '  +1.28%  '
$ws.Range('D9').Value = '''0.07244'
$ws.Range('E9').Value = '  +0.78%  '
$ws.Range('D10').Value = '''21.15'
$ws.Range('E10').Value = '  +2.15%  '
$ws.Range('D11').Value = '''0.8984'
$ws.Range('E11').Value = '  -0.01%  '
$ws.Range('D12').Value = '''0.08391'
$ws.Range('E12').Value = '  +11.06%  '
$ws.Range('D13').Value = '1.910.20'
$ws.Range('E13').Value = '  +1.98%  '
$ws.Range('D14').Value = '''94.55'
$ws.Range('E14').Value = '  -0.34%  '
$ws.Range('D15').Value = '''5.265'
$ws.Range('E15').Value = '  +0.19%  '
$ws.Range('D16').Value = '''1.001'
$ws.Range('E16').Value = '  +0.14%  '
$ws.Range('D17').Value = '''0.000008587'
$ws.Range('E17').Value = '  +1.16%  '
$ws.Range('E18').Value = '  +1.59%  '
$ws.Range('E19').Value = '  +0.11%  '
$ws.Range('D20').Value = '27.218.07'
$ws.Range('E20').Value = '  +0.93%  '
$ws.Range('D21').Value = '''5.053'
$ws.Range('E21').Value = '  +0.49%  '
$ws.Range('D22').Value = '2.147.84'
$ws.Range('E22').Value = '  +2.17%  '
$ws.Range('E23').Value = '  +1.69%  '
$ws.Range('D24').Value = '''6.412'
$ws.Range('E24').Value = '  -0.24%  '
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').Value = '''146.61'
$ws.Range('E25').Value = '  +0.39%  '
$ws.Range('B26').Value = 'LidoDAOToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D26').Value = '''2.279'
$ws.Range('E26').Value = '  +7.83%  '
$ws.Range('D27').Value = '''1.756'
$ws.Range('E27').Value = '  -1.45%  '
$ws.Range('D28').Value = '''18.10'
$ws.Range('E28').Value = '  +0.41%  '
$ws.Range('D29').Value = '''114.77'
$ws.Range('E29').Value = '  +0.12%  '
$ws.Range('D30').Value = '''4.919'
$ws.Range('E30').Value = '  +0.19%  '
$ws.Range('D31').Value = '''4.782'
$ws.Range('E31').Value = '  +0.58%  '
$ws.Range('D32').Value = '''0.09214'
$ws.Range('E32').Value = '  +0.43%  '
$ws.Range('D33').Value = '''0.8162'
$ws.Range('E33').Value = '  +8.42%  '
$ws.Range('D34').Value = '''0.05053'
$ws.Range('E34').Value = '  +0.51%  '
$ws.Range('D35').Value = '''1.235'
$ws.Range('E35').Value = '  +5.38%  '
$ws.Range('D36').Value = '''2.955'
$ws.Range('E36').Value = '  -1.19%  '
$ws.Range('D37').Value = '''3.382'
$ws.Range('E37').Value = '  +3.47%  '
$ws.Range('D38').Value = '''2.561'
$ws.Range('E38').Value = '  +2.95%  '
$ws.Range('D39').Value = '''0.5686'
$ws.Range('E39').Value = '  +2.08%  '
$ws.Range('D40').Value = '''0.01975'
$ws.Range('E40').Value = '  -0.93%  '
$ws.Range('D41').Value = '''1.073'
$ws.Range('E41').Value = '  +0.15%  '
$ws.Range('D42').Value = '''6.643'
$ws.Range('E42').Value = '  +1.02%  '
$ws.Range('D43').Value = '''8.926'
$ws.Range('E43').Value = '  +2.36%  '
$ws.Range('D44').Value = '''118.17'
$ws.Range('E44').Value = '  +2.00%  '
$ws.Range('D45').Value = '''0.1509'
$ws.Range('E45').Value = '  +0.77%  '
$ws.Range('D46').Value = '''0.4820'
$ws.Range('E46').Value = '  +1.05%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = '''10.18'
$ws.Range('E47').Value = '  +0.46%  '
$ws.Range('B48').Value = 'PaxDollar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D48').Value = '''1.001'
$ws.Range('E48').Value = '  +0.17%  '
$ws.Range('D49').Value = '''1.610'
$ws.Range('E49').Value = '  +3.06%  '
$ws.Range('D50').Value = '''37.40'
$ws.Range('E50').Value = '  +0.81%  '
$ws.Range('D51').Value = '''63.61'
$ws.Range('E51').Value = '  +0.31%  '
